# German translation pass for
# "Email 4-2 [TEMPLATE] Partner email – reminder to submit documents.docx"
#
# The template reuses the same English sentences more than once (e.g.
# "Don't forget to send your documents" and "We look forward to seeing
# you there!" each appear twice — once per copy of the email), and some
# look-alike English text must stay untouched (the "English" inside the
# language-picker hyperlink in paragraph 1, and the two longer
# "If you have any questions..." sentences that continue with extra
# wording). So rather than a document-wide Find/Replace, walk every
# paragraph and only replace it when its full (trimmed) text is an exact
# match for one of the known source strings — this naturally translates
# every occurrence that should change while leaving look-alikes alone.

$d = $word.ActiveDocument

$map = @{
    "English" = "Englisch";
    "Don’t forget to send your documents" = "Vergessen Sie nicht, Ihre Dokumente zu schicken";
    "If you have any questions, please contact your country manager." = "Wenn Sie Fragen haben, wenden Sie sich bitte an Ihren Ländermanager.";
    "We look forward to seeing you there!" = "Wir freuen uns darauf, Sie dort zu sehen!";
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $rng = $d.Paragraphs.Item($i).Range
    $txt = $rng.Text.TrimEnd([char]13, [char]7)
    if ($map.ContainsKey($txt)) {
        $newText = $map[$txt]
        $rng.Find.Execute($txt, $true, $true, $false, $false, $false, `
                           $true, 1, $false, $newText, 2)
    }
}
